$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1
$ws.Cells.Item(1, 1).Value = "ZoneNameUpdated"
$ws.Cells.Item(1, 2).Value = "ZoneDescriptionUpdated"
$ws.Cells.Item(1, 3).Value = "10/25/2022 8:21:16 PM"
$ws.Cells.Item(1, 5).Value = "https://connectedoffice-devicemanagement.azurewebsites.net/images/edit.png"
$ws.Cells.Item(1, 6).Value = "https://connectedoffice-devicemanagement.azurewebsites.net/Zones/Edit/160f93a1-5bf4-4b75-a93c-d774802f59f7"
$ws.Cells.Item(1, 8).Value = "https://connectedoffice-devicemanagement.azurewebsites.net/images/details.png"
$ws.Cells.Item(1, 9).Value = "https://connectedoffice-devicemanagement.azurewebsites.net/Zones/Details/160f93a1-5bf4-4b75-a93c-d774802f59f7"
$ws.Cells.Item(1, 11).Value = "https://connectedoffice-devicemanagement.azurewebsites.net/images/delete.png"
$ws.Cells.Item(1, 12).Value = "https://connectedoffice-devicemanagement.azurewebsites.net/Zones/Delete/160f93a1-5bf4-4b75-a93c-d774802f59f7"
$ws.Cells.Item(1, 13).Value = "https://connectedoffice-devicemanagement.azurewebsites.net/Devices/Details/6c8c8366-a8dc-4ed1-82a2-f2ef67cd97a5"
$ws.Cells.Item(1, 15).Value = "https://connectedoffice-devicemanagement.azurewebsites.net/images/delete.png"
$ws.Cells.Item(1, 16).Value = "https://connectedoffice-devicemanagement.azurewebsites.net/Devices/Delete/6c8c8366-a8dc-4ed1-82a2-f2ef67cd97a5"

# Row 2
$ws.Cells.Item(2, 1).Value = "Boilermaker RoomUpdated"
$ws.Cells.Item(2, 2).Value = "2-870 - Sculpture/OrnamentalUpdated"
$ws.Cells.Item(2, 3).Value = "10/25/2022 8:15:04 PM"
$ws.Cells.Item(2, 5).Value = "https://connectedoffice-devicemanagement.azurewebsites.net/images/edit.png"
$ws.Cells.Item(2, 6).Value = "https://connectedoffice-devicemanagement.azurewebsites.net/Zones/Edit/d554a358-e10f-4d98-8c93-62bd2d1c5d46"
$ws.Cells.Item(2, 8).Value = "https://connectedoffice-devicemanagement.azurewebsites.net/images/details.png"
$ws.Cells.Item(2, 9).Value = "https://connectedoffice-devicemanagement.azurewebsites.net/Zones/Details/d554a358-e10f-4d98-8c93-62bd2d1c5d46"
$ws.Cells.Item(2, 11).Value = "https://connectedoffice-devicemanagement.azurewebsites.net/images/delete.png"
$ws.Cells.Item(2, 12).Value = "https://connectedoffice-devicemanagement.azurewebsites.net/Zones/Delete/d554a358-e10f-4d98-8c93-62bd2d1c5d46"

# Row 3
$ws.Cells.Item(3, 1).Value = "ZoneNameUpdated"
$ws.Cells.Item(3, 2).Value = "ZoneDescriptionUpdated"
$ws.Cells.Item(3, 3).Value = "10/25/2022 8:14:53 PM"
$ws.Cells.Item(3, 5).Value = "https://connectedoffice-devicemanagement.azurewebsites.net/images/edit.png"
$ws.Cells.Item(3, 6).Value = "https://connectedoffice-devicemanagement.azurewebsites.net/Zones/Edit/efca46e9-6b5b-4212-95fe-83cc86bef3ba"
$ws.Cells.Item(3, 8).Value = "https://connectedoffice-devicemanagement.azurewebsites.net/images/details.png"
$ws.Cells.Item(3, 9).Value = "https://connectedoffice-devicemanagement.azurewebsites.net/Zones/Details/efca46e9-6b5b-4212-95fe-83cc86bef3ba"
$ws.Cells.Item(3, 11).Value = "https://connectedoffice-devicemanagement.azurewebsites.net/images/delete.png"
$ws.Cells.Item(3, 12).Value = "https://connectedoffice-devicemanagement.azurewebsites.net/Zones/Delete/efca46e9-6b5b-4212-95fe-83cc86bef3ba"

# Row 4
$ws.Cells.Item(4, 1).Value = "Safety OfficeUpdated"
$ws.Cells.Item(4, 2).Value = "17-030 - BondUpdated"
$ws.Cells.Item(4, 3).Value = "10/25/2022 7:20:45 PM"
$ws.Cells.Item(4, 5).Value = "https://connectedoffice-devicemanagement.azurewebsites.net/images/edit.png"
$ws.Cells.Item(4, 6).Value = "https://connectedoffice-devicemanagement.azurewebsites.net/Zones/Edit/4201f294-c200-4e75-b4ac-6a45151601c2"
$ws.Cells.Item(4, 8).Value = "https://connectedoffice-devicemanagement.azurewebsites.net/images/details.png"
$ws.Cells.Item(4, 9).Value = "https://connectedoffice-devicemanagement.azurewebsites.net/Zones/Details/4201f294-c200-4e75-b4ac-6a45151601c2"
$ws.Cells.Item(4, 11).Value = "https://connectedoffice-devicemanagement.azurewebsites.net/images/delete.png"
$ws.Cells.Item(4, 12).Value = "https://connectedoffice-devicemanagement.azurewebsites.net/Zones/Delete/4201f294-c200-4e75-b4ac-6a45151601c2"

# Row 5
$ws.Cells.Item(5, 1).Value = "Tile Setting BenchUpdated"
$ws.Cells.Item(5, 2).Value = "1-570 - Temporary ControlsUpdated"
$ws.Cells.Item(5, 3).Value = "10/25/2022 7:21:11 PM"
$ws.Cells.Item(5, 5).Value = "https://connectedoffice-devicemanagement.azurewebsites.net/images/edit.png"
$ws.Cells.Item(5, 6).Value = "https://connectedoffice-devicemanagement.azurewebsites.net/Zones/Edit/d62a79c8-7968-40d0-bee9-943c0983752e"
$ws.Cells.Item(5, 8).Value = "https://connectedoffice-devicemanagement.azurewebsites.net/images/details.png"
$ws.Cells.Item(5, 9).Value = "https://connectedoffice-devicemanagement.azurewebsites.net/Zones/Details/d62a79c8-7968-40d0-bee9-943c0983752e"
$ws.Cells.Item(5, 11).Value = "https://connectedoffice-devicemanagement.azurewebsites.net/images/delete.png"
$ws.Cells.Item(5, 12).Value = "https://connectedoffice-devicemanagement.azurewebsites.net/Zones/Delete/d62a79c8-7968-40d0-bee9-943c0983752e"

# Row 6
$ws.Cells.Item(6, 1).Value = "Boilermaker RoomUpdated"
$ws.Cells.Item(6, 2).Value = "2-870 - Sculpture/OrnamentalUpdated"
$ws.Cells.Item(6, 3).Value = "10/25/2022 7:20:40 PM"
$ws.Cells.Item(6, 5).Value = "https://connectedoffice-devicemanagement.azurewebsites.net/images/edit.png"
$ws.Cells.Item(6, 6).Value = "https://connectedoffice-devicemanagement.azurewebsites.net/Zones/Edit/48dc1bd2-caf3-4282-a254-a674fafc42ec"
$ws.Cells.Item(6, 8).Value = "https://connectedoffice-devicemanagement.azurewebsites.net/images/details.png"
$ws.Cells.Item(6, 9).Value = "https://connectedoffice-devicemanagement.azurewebsites.net/Zones/Details/48dc1bd2-caf3-4282-a254-a674fafc42ec"
$ws.Cells.Item(6, 11).Value = "https://connectedoffice-devicemanagement.azurewebsites.net/images/delete.png"
$ws.Cells.Item(6, 12).Value = "https://connectedoffice-devicemanagement.azurewebsites.net/Zones/Delete/48dc1bd2-caf3-4282-a254-a674fafc42ec"

# Row 7
$ws.Cells.Item(7, 1).Value = "Labor OfficeUpdated"
$ws.Cells.Item(7, 2).Value = "1-523 - Sanitary FacilitiesUpdated"
$ws.Cells.Item(7, 3).Value = "10/25/2022 7:21:01 PM"
$ws.Cells.Item(7, 5).Value = "https://connectedoffice-devicemanagement.azurewebsites.net/images/edit.png"
$ws.Cells.Item(7, 6).Value = "https://connectedoffice-devicemanagement.azurewebsites.net/Zones/Edit/0dca1b39-6e01-423f-a734-c56f12ca7b53"
$ws.Cells.Item(7, 8).Value = "https://connectedoffice-devicemanagement.azurewebsites.net/images/details.png"
$ws.Cells.Item(7, 9).Value = "https://connectedoffice-devicemanagement.azurewebsites.net/Zones/Details/0dca1b39-6e01-423f-a734-c56f12ca7b53"
$ws.Cells.Item(7, 11).Value = "https://connectedoffice-devicemanagement.azurewebsites.net/images/delete.png"
$ws.Cells.Item(7, 12).Value = "https://connectedoffice-devicemanagement.azurewebsites.net/Zones/Delete/0dca1b39-6e01-423f-a734-c56f12ca7b53"

# Row 8
$ws.Cells.Item(8, 1).Value = "Stucco Mason BuildingUpdated"
$ws.Cells.Item(8, 2).Value = "2-750 - Concrete Pads and WalksUpdated"
$ws.Cells.Item(8, 3).Value = "10/25/2022 7:20:51 PM"
$ws.Cells.Item(8, 5).Value = "https://connectedoffice-devicemanagement.azurewebsites.net/images/edit.png"
$ws.Cells.Item(8, 6).Value = "https://connectedoffice-devicemanagement.azurewebsites.net/Zones/Edit/1122805a-abbc-4fc5-85c1-ef88b3cd4e2a"
$ws.Cells.Item(8, 8).Value = "https://connectedoffice-devicemanagement.azurewebsites.net/images/details.png"
$ws.Cells.Item(8, 9).Value = "https://connectedoffice-devicemanagement.azurewebsites.net/Zones/Details/1122805a-abbc-4fc5-85c1-ef88b3cd4e2a"
$ws.Cells.Item(8, 11).Value = "https://connectedoffice-devicemanagement.azurewebsites.net/images/delete.png"
$ws.Cells.Item(8, 12).Value = "https://connectedoffice-devicemanagement.azurewebsites.net/Zones/Delete/1122805a-abbc-4fc5-85c1-ef88b3cd4e2a"

# Row 9
$ws.Cells.Item(9, 1).Value = "ZoneNameUpdated"
$ws.Cells.Item(9, 2).Value = "ZoneDescriptionUpdated"
$ws.Cells.Item(9, 3).Value = "10/25/2022 7:20:30 PM"
$ws.Cells.Item(9, 5).Value = "https://connectedoffice-devicemanagement.azurewebsites.net/images/edit.png"
$ws.Cells.Item(9, 6).Value = "https://connectedoffice-devicemanagement.azurewebsites.net/Zones/Edit/1c39bed7-abd1-4430-9551-ff464a872b2b"
$ws.Cells.Item(9, 8).Value = "https://connectedoffice-devicemanagement.azurewebsites.net/images/details.png"
$ws.Cells.Item(9, 9).Value = "https://connectedoffice-devicemanagement.azurewebsites.net/Zones/Details/1c39bed7-abd1-4430-9551-ff464a872b2b"
$ws.Cells.Item(9, 11).Value = "https://connectedoffice-devicemanagement.azurewebsites.net/images/delete.png"
$ws.Cells.Item(9, 12).Value = "https://connectedoffice-devicemanagement.azurewebsites.net/Zones/Delete/1c39bed7-abd1-4430-9551-ff464a872b2b"
